# Refresh the "feature importance" table on Sheet1 with the finished
# Random Forest model's results, re-ranked by impDecrease (desc.), and
# append one more row for a manually-added feature ("rat_spcsrc") that
# was kept out of the RF run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old data rows first. Re-entering everything from a clean
# range (rather than overwriting cell-by-cell) keeps the shared-string
# table free of now-unused entries.
$ws.Range("A2:D32").ClearContents()

# Row 33 is brand new -- clone the formatting of the last existing data
# row (bold/centered rank style in column A) before filling it in.
$ws.Range("A32").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$rows = @(
    @(0,  "sstk_std",             0.05152090202630415, "Random Forest"),
    @(1,  "teq",                  0.04210914267773873, "Random Forest"),
    @(2,  "seq",                  0.04113480038972558, "Random Forest"),
    @(3,  "rest_sum_diff",        0.02964086773085264, "Random Forest"),
    @(4,  "xopr",                 0.02894210148331851, "Random Forest"),
    @(5,  "ceq",                  0.02673327775011868, "Random Forest"),
    @(6,  "caps",                 0.02374550671631139, "Random Forest"),
    @(7,  "tstk",                 0.02261428377865916, "Random Forest"),
    @(8,  "fopo_std",             0.02202374152521238, "Random Forest"),
    @(9,  "sec_trt1m_std",        0.02085166940573271, "Random Forest"),
    @(10, "ceqt",                 0.02081067191182673, "Random Forest"),
    @(11, "icapt",                0.02018532661059256, "Random Forest"),
    @(12, "rest_count_of_diffs",  0.02006966905846362, "Random Forest"),
    @(13, "lse",                  0.01844682961883291, "Random Forest"),
    @(14, "pi_std",               0.0179690776523493,  "Random Forest"),
    @(15, "revt",                 0.01768615437752864, "Random Forest"),
    @(16, "sstk",                 0.01747198953836201, "Random Forest"),
    @(17, "st_per_growth",        0.01659163279284432, "Random Forest"),
    @(18, "at",                   0.01653665368324137, "Random Forest"),
    @(19, "rest_count",           0.01653338496751571, "Random Forest"),
    @(20, "cogs",                 0.01642483665464855, "Random Forest"),
    @(21, "gp",                   0.01639626005398251, "Random Forest"),
    @(22, "xsga",                 0.01488677887659137, "Random Forest"),
    @(23, "rect_std",             0.01412363442861623, "Random Forest"),
    @(24, "cogs_std",             0.01403655662386116, "Random Forest"),
    @(25, "lct",                  0.01353114676459853, "Random Forest"),
    @(26, "dilavx_std",           0.01110092391154593, "Random Forest"),
    @(27, "tstk_std",             0.01082068248664309, "Random Forest"),
    @(28, "auop",                 0.01058660278150156, "Random Forest"),
    @(29, "dlc",                  0.01045068502105618, "Random Forest"),
    @(30, "xsga_std",             0.01028382202854827, "Random Forest"),
    @(31, "rat_spcsrc",           0.001111,             "Manual Addition")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
